# Devdutt Padikkal.xlsx - append the batting-innings table rows 17-31.
# These are a re-ordered repeat of the existing rows 2-16 (same stats,
# appended again further down the sheet), so we copy each destination row
# from its corresponding existing source row rather than retyping values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking columns (G:K) of the new rows as Text so
# that values such as "50.00" / "150.00" keep their original text
# representation instead of being auto-coerced into numbers (which would
# silently drop the trailing zeros). Columns A:F never look numeric, so
# they don't need this and are left on the sheet's default format.
$ws.Range("G17:K31").NumberFormat = "@"

# row 17 <- row 3
$ws.Range("A17:F17").Value = $ws.Range("A3:F3").Value2
$ws.Range("G17:K17").Value = $ws.Range("G3:K3").Value2

# row 18 <- row 4
$ws.Range("A18:F18").Value = $ws.Range("A4:F4").Value2
$ws.Range("G18:K18").Value = $ws.Range("G4:K4").Value2

# row 19 <- row 8
$ws.Range("A19:F19").Value = $ws.Range("A8:F8").Value2
$ws.Range("G19:K19").Value = $ws.Range("G8:K8").Value2

# row 20 <- row 5
$ws.Range("A20:F20").Value = $ws.Range("A5:F5").Value2
$ws.Range("G20:K20").Value = $ws.Range("G5:K5").Value2

# row 21 <- row 11
$ws.Range("A21:F21").Value = $ws.Range("A11:F11").Value2
$ws.Range("G21:K21").Value = $ws.Range("G11:K11").Value2

# row 22 <- row 6
$ws.Range("A22:F22").Value = $ws.Range("A6:F6").Value2
$ws.Range("G22:K22").Value = $ws.Range("G6:K6").Value2

# row 23 <- row 16
$ws.Range("A23:F23").Value = $ws.Range("A16:F16").Value2
$ws.Range("G23:K23").Value = $ws.Range("G16:K16").Value2

# row 24 <- row 13
$ws.Range("A24:F24").Value = $ws.Range("A13:F13").Value2
$ws.Range("G24:K24").Value = $ws.Range("G13:K13").Value2

# row 25 <- row 14
$ws.Range("A25:F25").Value = $ws.Range("A14:F14").Value2
$ws.Range("G25:K25").Value = $ws.Range("G14:K14").Value2

# row 26 <- row 15
$ws.Range("A26:F26").Value = $ws.Range("A15:F15").Value2
$ws.Range("G26:K26").Value = $ws.Range("G15:K15").Value2

# row 27 <- row 7
$ws.Range("A27:F27").Value = $ws.Range("A7:F7").Value2
$ws.Range("G27:K27").Value = $ws.Range("G7:K7").Value2

# row 28 <- row 10
$ws.Range("A28:F28").Value = $ws.Range("A10:F10").Value2
$ws.Range("G28:K28").Value = $ws.Range("G10:K10").Value2

# row 29 <- row 9
$ws.Range("A29:F29").Value = $ws.Range("A9:F9").Value2
$ws.Range("G29:K29").Value = $ws.Range("G9:K9").Value2

# row 30 <- row 12
$ws.Range("A30:F30").Value = $ws.Range("A12:F12").Value2
$ws.Range("G30:K30").Value = $ws.Range("G12:K12").Value2

# row 31 <- row 2
$ws.Range("A31:F31").Value = $ws.Range("A2:F2").Value2
$ws.Range("G31:K31").Value = $ws.Range("G2:K2").Value2
